$d = $word.ActiveDocument

# --- Step 1: paragraph "Podstanar dobija formu u koju unosi naziv i opis obaveštenja" ---
# Trim it down to "...naziv obaveštenja" and add a new bullet right after it:
# "Podstanar unosi opis obaveštenja" (same list style as the paragraph it follows).

$r1 = $d.Content
$r1.Find.Execute("naziv i opis obaveštenja", $true, $false, $false, $false, $false, $true, 0, $false, "naziv obaveštenja", 2) | Out-Null

$r2 = $d.Content
$found2 = $r2.Find.Execute("naziv obaveštenja", $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
if ($found2) {
    $r2.InsertAfter("`rPodstanar unosi opis obaveštenja")
}

# --- Step 2: paragraph "Pritiskom na dugme "Okači obaveštenje", poruka se kači na oglasnu           tablu" ---
# Update the button caption, tighten the long run of spaces, then add a new bullet right after it:
# "Podstanar se vraća u korak 2.2.1.1. i omogućen mu je ponovni unos obaveštenja"

$r3 = $d.Content
$r3.Find.Execute("“Okači obaveštenje”", $true, $false, $false, $false, $false, $true, 0, $false, "“Potvrdi kačenje obaveštenja na oglasnu tablu”", 2) | Out-Null

$r4 = $d.Content
$r4.Find.Execute("oglasnu           tablu", $true, $false, $false, $false, $false, $true, 0, $false, "oglasnu   tablu", 2) | Out-Null

$r5 = $d.Content
$found5 = $r5.Find.Execute("oglasnu   tablu", $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
if ($found5) {
    $r5.InsertAfter("`rPodstanar se vraća u korak 2.2.1.1. i omogućen mu je ponovni unos obaveštenja")
}

# Re-create the "_GoBack" bookmark (Word keeps exactly one, tracking the last edit)
# at the end of the text we just inserted.
$r6 = $d.Content
$found6 = $r6.Find.Execute("ponovni unos obaveštenja", $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
if ($found6) {
    $r6.Collapse(0)
    $d.Bookmarks.Add("_GoBack", $r6) | Out-Null
}
